$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename precondition "Guest" -> "Guest user" for rows 4-22
foreach ($r in 4..22) {
    $ws.Range("D$r").Value = "Guest user"
}

# Fill in Test Steps (E) / Expected Result (F) for the six Main Navigation Menu
# category rows that previously had none: Makeup(9), Skincare(10), Fragrance(11),
# Men(12), Hair Care(13), Books(14).

$ws.Range("E9").Value = "1. Open https://abantecart.codifyme.co.nz`n2. Check if MAKEUP Menu link is working`n3. Mouseover MAKEUP Menu`n4. Check if all links in the MAKEUP Menu dropdown are working"
$ws.Range("F9").Value = "2. MAKEUP Menu button link should be working.`n3. Dropdown appears .`n4. All dropdown links should be working."

$ws.Range("E10").Value = "1. Open https://abantecart.codifyme.co.nz`n2. Check if SKINCARE Menu link is working`n3. Mouseover SKINCARE Menu`n4. Check if all links in the SKINCARE Menu dropdown are working"
$ws.Range("F10").Value = "2. SKINCARE Menu button link should be working.`n3. Dropdown appears .`n4. All dropdown links should be working."

$ws.Range("E11").Value = "1. Open https://abantecart.codifyme.co.nz`n2. Check if FRAGRANCE Menu link is working`n3. Mouseover FRAGRANCE Menu`n4. Check if all links in the FRAGRANCE Menu dropdown are working"
$ws.Range("F11").Value = "2. FRAGRANCE Menu button link should be working.`n3. Dropdown appears .`n4. All dropdown links should be working."

$ws.Range("E12").Value = "1. Open https://abantecart.codifyme.co.nz`n2. Check if MEN Menu link is working`n3. Mouseover MEN Menu`n4. Check if all links in the MEN Menu dropdown are working"
$ws.Range("F12").Value = "2. MEN Menu button link should be working.`n3. Dropdown appears .`n4. All dropdown links should be working."

$ws.Range("E13").Value = "1. Open https://abantecart.codifyme.co.nz`n2. Check if HAIR CARE Menu link is working`n3. Mouseover HAIR CARE Menu`n4. Check if all links in the HAIR CARE Menu dropdown are working"
$ws.Range("F13").Value = "2. HAIR CARE Menu button link should be working.`n3. Dropdown appears .`n4. All dropdown links should be working."

$ws.Range("E14").Value = "1. Open https://abantecart.codifyme.co.nz`n2. Check if BOOKS Menu link is working`n3. Mouseover BOOKS Menu`n4. Check if all links in the BOOKS Menu dropdown are working"
$ws.Range("F14").Value = "2. BOOKS Menu button link should be working.`n3. Dropdown appears .`n4. All dropdown links should be working."

# Row heights for the newly-populated rows (9-14) match the other wrapped-text rows (75pt)
foreach ($r in 9..14) {
    $ws.Rows.Item($r).RowHeight = 75
}

# Frozen pane top-left cell moves from A6 to A8 to reflect the newly inserted rows
$ws.Application.ActiveWindow.FreezePanes = $false
$ws.Range("A8").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
